$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Title'
$ws.Range("B1").Value = 'Content'
$ws.Range("A2").Value = 'Tesla Stock Is Down -- Regulators Are Taking Aim at ''Full Self-Driving'' - Yahoo Finance'
$ws.Range("B2").Value = '- Best non-toxic cutting boards. - Best high-yield savings accounts. ### New on Yahoo. # Yahoo Finance. Tesla (NASDAQ:TSLA) stock is down about 0.5% on Wednesday morning as the company faces a regulatory headache in California after the Department of Motor Vehicles said the automaker used misleading terms such as Autopilot and Full Self-Driving for its driver-assist system. The regulator gave Tesla 60 days to change marketing language or risk suspension of its dealer license, while sparing manufacturing permits for now. The DMV emphasized Tesla''s system remains Level 2, drivers must stay attentive, even as the company tests robo-taxi services with safety monitors aboard. Tesla began a limited self-driving cab service in Austin and continues tests elsewhere, keeping a safety monitor in the passenger seat for now. Analysts said simple rebranding could address the DMV''s concern, but noted ongoing safety and regulatory scrutiny could slow broader robotaxi plans. CEO Elon Musk has said Tesla aims to remove safety monitors by year-end, a timeline analysts call ambitious given regulatory hurdles.'
$ws.Range("A3").Value = 'Tesla Stock Is Down -- Regulators Are Taking Aim at ''Full Self-Driving'' - Yahoo! Finance Canada'
$ws.Range("B3").Value = 'About our ads  Advertising  Jobs. #### MARKETS LIVE BLOG. DOW, S&P 500, NASDAQ RISE AS WALL ST. WEIGHS JOBS DATA SIGNALS, ORACLE SINKS. Tesla (NASDAQ:TSLA) stock is down about 0.5% on Wednesday morning as the company faces a regulatory headache in California after the Department of Motor Vehicles said the automaker used misleading terms such as Autopilot and Full Self-Driving for its driver-assist system. The regulator gave Tesla 60 days to change marketing language or risk suspension of its dealer license, while sparing manufacturing permits for now. Test your thesis with our free DCF calculator. The DMV emphasized Tesla''s system remains Level 2, drivers must stay attentive, even as the company tests robo-taxi services with safety monitors aboard. Analysts said simple rebranding could address the DMV''s concern, but noted ongoing safety and regulatory scrutiny could slow broader robotaxi plans. CEO Elon Musk has said Tesla aims to remove safety monitors by year-end, a timeline analysts call ambitious given regulatory hurdles.'
$ws.Range("A4").Value = 'Tesla Stock (TSLA) Near Record High as Robotaxi “No-Occupant” Testing Ignites Fresh Bull Case Forecasts - ts2.tech'
$ws.Range("B4").Value = '# Tesla Stock (TSLA) Near Record High as Robotaxi “No-Occupant” Testing Ignites Fresh Bull Case Forecasts. (NASDAQ: TSLA) is back in the spotlight after a burst of autonomy-related headlines helped push the stock to the edge of record territory, reviving the long-running debate over whether investors are valuing Tesla primarily as an automaker — or as an AI-and-robotics platform. Tesla stock’s latest surge is a textbook example of how TSLA trades when autonomy momentum appears to accelerate: **the market rapidly re-prices the long-term robotaxi opportunity**, while skeptics emphasize the absence of widely published safety metrics, competitive pressure from Waymo, and the reality that Tesla’s near-term financials still depend on vehicle sales. 2. Tesla Stock (TSLA) Near Record High as Robotaxi “No-Occupant” Testing Ignites Fresh Bull Case Forecasts. ### Tesla Stock (TSLA) Near Record High as Robotaxi “No-Occupant” Testing Ignites Fresh Bull Case Forecasts.'
$ws.Range("A5").Value = 'Tesla Stock (TSLA) News on Dec. 14, 2025: Robotaxi “Three-Week” Countdown Collides With Slowing Sales—What Analysts Forecast Next - ts2.tech'
$ws.Range("B5").Value = 'Analyst Forecasts · Electric Vehicles · NASDAQ:TSLA · Robotaxi · Stock Market · Tesla. # Tesla Stock (TSLA) News on Dec. 14, 2025: Robotaxi “Three-Week” Countdown Collides With Slowing Sales—What Analysts Forecast Next. Tesla stock is holding up near $459 as robotaxi ambitions fuel optimism—but fresh U.S., Europe, and China sales data point to softer demand. The most market-moving headline in recent days came from Reuters, citing Cox Automotive estimates: Tesla’s U.S. sales fell nearly **23% year-over-year** in **November 2025** to about **39,800 vehicles**, the lowest monthly total since **January 2022**—even after Tesla rolled out cheaper “Standard” versions of the Model Y and Model 3. Cox estimated overall U.S. EV sales dropped **more than 41%** in November, lifting Tesla’s market share to **56.7%** from **43.1%** a year earlier. 3. Tesla Stock (TSLA) News on Dec. 14, 2025: Robotaxi “Three-Week” Countdown Collides With Slowing Sales—What Analysts Forecast Next. ### Tesla Stock (TSLA) News on Dec. 14, 2025: Robotaxi “Three-Week” Countdown Collides With Slowing Sales—What Analysts Forecast Next.'
$ws.Range("A6").Value = 'Tesla Stock (TSLA) Before US Market Open 15 December 2025: Latest News, Analyst Forecasts, and Key Catalysts to Watch - ts2.tech'
$ws.Range("B6").Value = '# Tesla Stock (TSLA) Before US Market Open 15 December 2025: Latest News, Analyst Forecasts, and Key Catalysts to Watch. Tesla stock heads into Monday’s US session with a familiar setup: strong share momentum and a sky-high valuation narrative tied to autonomy and robotics—colliding with fresh evidence of weakening demand in several core EV markets and heightened regulatory scrutiny around Full Self-Driving (FSD). Below is what to know about TSLA before the US stock market opens on **15/12/2025**, based on the most recent reporting, regulatory documents, and analyst commentary. This is a major “headline risk” factor for Monday and beyond: any indication the probe is escalating—or that Tesla is changing how it markets, gates, or deploys FSD features—can move the stock quickly. 2. Tesla Stock (TSLA) Before US Market Open 15 December 2025: Latest News, Analyst Forecasts, and Key Catalysts to Watch. ### Tesla Stock (TSLA) Before US Market Open 15 December 2025: Latest News, Analyst Forecasts, and Key Catalysts to Watch.'
